$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Paragraphs.Item(1).Range.Text = "2024-04-14 Sunday"

# Update each table cell value (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "81-3="
$t.Cell(1, 2).Range.Text = "37+48="
$t.Cell(1, 3).Range.Text = "29+53="
$t.Cell(1, 4).Range.Text = "9+27="
$t.Cell(1, 5).Range.Text = "91-83="

$t.Cell(2, 1).Range.Text = "89+5="
$t.Cell(2, 2).Range.Text = "68+6="
$t.Cell(2, 3).Range.Text = "44-29="
$t.Cell(2, 4).Range.Text = "32-7="
$t.Cell(2, 5).Range.Text = "93-19="

$t.Cell(3, 1).Range.Text = "81-2="
$t.Cell(3, 2).Range.Text = "86-47="
$t.Cell(3, 3).Range.Text = "66-59="
$t.Cell(3, 4).Range.Text = "44-16="
$t.Cell(3, 5).Range.Text = "9+35="

$t.Cell(4, 1).Range.Text = "80-34="
$t.Cell(4, 2).Range.Text = "37+9="
$t.Cell(4, 3).Range.Text = "38+15="
$t.Cell(4, 4).Range.Text = "10-7="
$t.Cell(4, 5).Range.Text = "28+37="

$t.Cell(5, 1).Range.Text = "95-26="
$t.Cell(5, 2).Range.Text = "76-19="
$t.Cell(5, 3).Range.Text = "4+29="
$t.Cell(5, 4).Range.Text = "94-75="
$t.Cell(5, 5).Range.Text = "56+26="

$t.Cell(6, 1).Range.Text = "39+24="
$t.Cell(6, 2).Range.Text = "52-43="
$t.Cell(6, 3).Range.Text = "38+57="
$t.Cell(6, 4).Range.Text = "27+25="
$t.Cell(6, 5).Range.Text = "46+49="

$t.Cell(7, 1).Range.Text = "37+7="
$t.Cell(7, 2).Range.Text = "29+5="
$t.Cell(7, 3).Range.Text = "30-4="
$t.Cell(7, 4).Range.Text = "8+3="
$t.Cell(7, 5).Range.Text = "19+38="

$t.Cell(8, 1).Range.Text = "17+48="
$t.Cell(8, 2).Range.Text = "24-17="
$t.Cell(8, 3).Range.Text = "19+24="
$t.Cell(8, 4).Range.Text = "66-47="
$t.Cell(8, 5).Range.Text = "16+16="

$t.Cell(9, 1).Range.Text = "67-28="
$t.Cell(9, 2).Range.Text = "19+33="
$t.Cell(9, 3).Range.Text = "3+28="
$t.Cell(9, 4).Range.Text = "14+27="
$t.Cell(9, 5).Range.Text = "75-29="

$t.Cell(10, 1).Range.Text = "8+75="
$t.Cell(10, 2).Range.Text = "86-77="
$t.Cell(10, 3).Range.Text = "74+19="
$t.Cell(10, 4).Range.Text = "55+17="
$t.Cell(10, 5).Range.Text = "28+14="

$t.Cell(11, 1).Range.Text = "82-56="
$t.Cell(11, 2).Range.Text = "85-7="
$t.Cell(11, 3).Range.Text = "97-79="
$t.Cell(11, 4).Range.Text = "97-79="
$t.Cell(11, 5).Range.Text = "14+29="

$t.Cell(12, 1).Range.Text = "42-37="
$t.Cell(12, 2).Range.Text = "46+28="
$t.Cell(12, 3).Range.Text = "48+47="
$t.Cell(12, 4).Range.Text = "41-37="
$t.Cell(12, 5).Range.Text = "62-3="

$t.Cell(13, 1).Range.Text = "31-24="
$t.Cell(13, 2).Range.Text = "94-17="
$t.Cell(13, 3).Range.Text = "21-4="
$t.Cell(13, 4).Range.Text = "98-89="
$t.Cell(13, 5).Range.Text = "94-25="

$t.Cell(14, 1).Range.Text = "69+24="
$t.Cell(14, 2).Range.Text = "53-26="
$t.Cell(14, 3).Range.Text = "10-2="
$t.Cell(14, 4).Range.Text = "43-35="
$t.Cell(14, 5).Range.Text = "19+16="

$t.Cell(15, 1).Range.Text = "75-67="
$t.Cell(15, 2).Range.Text = "3+68="
$t.Cell(15, 3).Range.Text = "78-39="
$t.Cell(15, 4).Range.Text = "39+37="
$t.Cell(15, 5).Range.Text = "31-7="

$t.Cell(16, 1).Range.Text = "4+89="
$t.Cell(16, 2).Range.Text = "16+36="
$t.Cell(16, 3).Range.Text = "19+25="
$t.Cell(16, 4).Range.Text = "26+65="
$t.Cell(16, 5).Range.Text = "36+29="

$t.Cell(17, 1).Range.Text = "40-29="
$t.Cell(17, 2).Range.Text = "39+22="
$t.Cell(17, 3).Range.Text = "67+27="
$t.Cell(17, 4).Range.Text = "40-38="
$t.Cell(17, 5).Range.Text = "55+37="

$t.Cell(18, 1).Range.Text = "74-5="
$t.Cell(18, 2).Range.Text = "58+18="
$t.Cell(18, 3).Range.Text = "95-36="
$t.Cell(18, 4).Range.Text = "86+6="
$t.Cell(18, 5).Range.Text = "62-45="

$t.Cell(19, 1).Range.Text = "85-37="
$t.Cell(19, 2).Range.Text = "69+13="
$t.Cell(19, 3).Range.Text = "27+54="
$t.Cell(19, 4).Range.Text = "45-17="
$t.Cell(19, 5).Range.Text = "10-2="

$t.Cell(20, 1).Range.Text = "6+9="
$t.Cell(20, 2).Range.Text = "37-19="
$t.Cell(20, 3).Range.Text = "9+39="
$t.Cell(20, 4).Range.Text = "29+4="
$t.Cell(20, 5).Range.Text = "62-54="
